$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.877.94"
$ws.Range("E2").Value = "  +2.65%  "
$ws.Range("D3").Value = "'1.869.92"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").Value = "'313.30"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "'1.011"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "'0.4826"
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("D8").Value = "'0.3826"
$ws.Range("E8").Value = "  +3.62%  "
$ws.Range("D9").Value = "'0.07372"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").Value = "'0.9391"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "'21.02"
$ws.Range("E11").Value = "  +5.55%  "
$ws.Range("D12").Value = "'0.07815"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").Value = "'1.902.81"
$ws.Range("E13").Value = "  +3.63%  "
$ws.Range("D14").Value = "'5.499"
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").Value = "'6.609"
$ws.Range("E15").Value = "  +1.87%  "
$ws.Range("D16").Value = "'90.84"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "'0.000008882"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("D19").Value = "'1.010"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").Value = "'28.079.18"
$ws.Range("E20").Value = "  +3.36%  "
$ws.Range("D21").Value = "'14.81"
$ws.Range("E21").Value = "  +1.31%  "
$ws.Range("D22").Value = "'5.121"
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("D23").Value = "'2.140.30"
$ws.Range("E23").Value = "  +3.91%  "
$ws.Range("D24").Value = "'10.83"
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("D25").Value = "'1.941"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").Value = "'156.62"
$ws.Range("E26").Value = "  +2.28%  "
$ws.Range("D27").Value = "'18.58"
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("D28").Value = "'2.054"
$ws.Range("E28").Value = "  +3.04%  "
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("D30").Value = "'4.983"
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("D31").Value = "'0.08920"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D32").Value = "'3.334"
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("D33").Value = "'1.222"
$ws.Range("E33").Value = "  +3.22%  "
$ws.Range("D34").Value = "'0.7660"
$ws.Range("E34").Value = "  +3.58%  "
$ws.Range("D35").Value = "'4.662"
$ws.Range("E35").Value = "  +2.96%  "
$ws.Range("D36").Value = "'2.723"
$ws.Range("E36").Value = "  +1.59%  "
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("D38").Value = "'0.02043"
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("D39").Value = "'0.5646"
$ws.Range("D40").Value = "'0.05370"
$ws.Range("E40").Value = "  +1.77%  "
$ws.Range("D41").Value = "'2.990"
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("D42").Value = "'7.062"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "'8.565"
$ws.Range("E43").Value = "  +3.02%  "
$ws.Range("D44").Value = "'0.1535"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("E45").Value = "  +1.97%  "
$ws.Range("D46").Value = "'0.4894"
$ws.Range("E46").Value = "  +3.16%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'105.42"
$ws.Range("E47").Value = "  +3.24%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "'1.011"
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("D49").Value = "'1.675"
$ws.Range("E49").Value = "  +3.33%  "
$ws.Range("D50").Value = "'67.71"
$ws.Range("E50").Value = "  +2.77%  "
$ws.Range("D51").Value = "'0.06108"
$ws.Range("E51").Value = "  +0.69%  "
